$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3.682053159252038
$ws.Range("F2").Value = 5.41

$ws.Range("E3").Value = 5.965051755218097
$ws.Range("F3").Value = 6.385616604281449

$ws.Range("E4").Value = 7.115699172182123
$ws.Range("F4").Value = 6.705381561317198

$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 15

$ws.Range("F6").Value = 7.313857748729716

$ws.Range("E7").Value = 5.306426188541749
$ws.Range("F7").Value = 5.981500086660821

$ws.Range("F8").Value = 3.56

$ws.Range("E9").Value = 1.981494588585559
$ws.Range("F9").Value = 1.104760956810323
